$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Site Design"
$ws.Range("A9").Value = "http://data.baltimoresun.com/ga-session-issues-2015/"
$ws.Range("B9").Value = "From http://damihere.com/"
$ws.Range("A10").Value = "http://interactives.dallasnews.com/2015/topletz-homes/"

$ws.Columns.Item(2).ColumnWidth = 92.5

$ws.Range("A13").Select()
